$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
  $cell = $ws.Cells.Item($row, $col)
  $cell.Value = "'" + $val
  $cell.Style = "Normal"
}

# Update Price (D) and Volume(1h) (E) columns
Set-TextCell 2 4 '46.891.31'
Set-TextCell 2 5 '  +5.59%  '
Set-TextCell 3 4 '2.307.30'
Set-TextCell 3 5 '  +3.53%  '
Set-TextCell 4 5 '  -0.66%  '
Set-TextCell 5 4 '303.94'
Set-TextCell 5 5 '  +1.45%  '
Set-TextCell 6 4 '101.81'
Set-TextCell 6 5 '  +12.24%  '
Set-TextCell 7 5 '  +1.27%  '
Set-TextCell 8 4 '0.999'
Set-TextCell 8 5 '  -0.54%  '
Set-TextCell 9 4 '0.531'
Set-TextCell 9 5 '  +7.51%  '
Set-TextCell 10 4 '36.93'
Set-TextCell 10 5 '  +11.40%  '
Set-TextCell 11 4 '0.0805'
Set-TextCell 11 5 '  +2.29%  '
Set-TextCell 12 4 '7.49'
Set-TextCell 12 5 '  +7.21%  '
Set-TextCell 13 5 '  +0.34%  '
Set-TextCell 14 4 '2.652.69'
Set-TextCell 14 5 '  +3.26%  '
Set-TextCell 15 4 '2.303.20'
Set-TextCell 15 5 '  +3.42%  '
Set-TextCell 16 4 '14.04'
Set-TextCell 16 5 '  +4.12%  '
Set-TextCell 17 4 '0.823'
Set-TextCell 17 5 '  +5.23%  '
Set-TextCell 18 4 '46.847.22'
Set-TextCell 18 5 '  +5.61%  '
Set-TextCell 19 4 '13.68'
Set-TextCell 19 5 '  +23.41%  '
Set-TextCell 20 4 '0.0₃0952'
Set-TextCell 20 5 '  +4.20%  '
Set-TextCell 21 4 '6.13'
Set-TextCell 21 5 '  +3.20%  '
Set-TextCell 22 4 '66.99'
Set-TextCell 22 5 '  +3.44%  '
Set-TextCell 23 4 '249.27'
Set-TextCell 23 5 '  +3.42%  '
Set-TextCell 24 5 '  +5.11%  '
Set-TextCell 25 5 '  +5.32%  '
Set-TextCell 26 4 '0.998'
Set-TextCell 26 5 '  -1.35%  '
Set-TextCell 27 4 '44.08'
Set-TextCell 27 5 '  +14.72%  '
Set-TextCell 28 5 '  +1.61%  '
Set-TextCell 29 4 '9.95'
Set-TextCell 29 5 '  +6.01%  '
Set-TextCell 30 4 '20.25'
Set-TextCell 30 5 '  +3.57%  '
Set-TextCell 31 4 '5.80'
Set-TextCell 31 5 '  +7.24%  '
Set-TextCell 32 4 '147.36'
Set-TextCell 32 5 '  -0.85%  '
Set-TextCell 33 4 '0.0802'
Set-TextCell 33 5 '  +6.50%  '
Set-TextCell 34 5 '  +3.18%  '
Set-TextCell 35 4 '3.19'
Set-TextCell 35 5 '  +11.88%  '
Set-TextCell 36 5 '  +10.88%  '
Set-TextCell 37 5 '  +3.08%  '
Set-TextCell 38 4 '1.81'
Set-TextCell 38 5 '  +6.54%  '
Set-TextCell 39 4 '16.11'
Set-TextCell 39 5 '  +21.46%  '
Set-TextCell 40 5 '  +13.33%  '
Set-TextCell 41 5 '  +8.91%  '
Set-TextCell 42 4 '0.0307'
Set-TextCell 42 5 '  +0.68%  '
Set-TextCell 46 4 '88.60'
Set-TextCell 46 5 '  +18.64%  '
Set-TextCell 47 4 '0.198'
Set-TextCell 47 5 '  +9.81%  '
Set-TextCell 48 4 '74.99'
Set-TextCell 48 5 '  +11.92%  '
Set-TextCell 49 4 '4.91'
Set-TextCell 49 5 '  +9.27%  '
Set-TextCell 50 4 '97.30'
Set-TextCell 50 5 '  +2.95%  '
Set-TextCell 51 4 '8.08'
Set-TextCell 51 5 '  +6.09%  '

# Rows 43-45 rotate: Stacks -> row43, FirstDigitalUSD -> row44, Maker -> row45
Set-TextCell 43 2 'Stacks'
Set-TextCell 43 3 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell 43 4 '2.00'
Set-TextCell 43 5 '  +12.28%  '
Set-TextCell 44 2 'FirstDigitalUSD'
Set-TextCell 44 3 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell 44 4 '0.998'
Set-TextCell 44 5 '  -0.73%  '
Set-TextCell 45 2 'Maker'
Set-TextCell 45 3 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 45 4 '1.865.80'
Set-TextCell 45 5 '  +2.56%  '
